$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential note date from 2021-03-19 to 2021-03-22
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-22 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2450297376388265
$ws.Range("E2").Value = 0.01438210227272729

$ws.Range("D3").Value = 0.4952977506214857
$ws.Range("E3").Value = 0.00393258426966292

$ws.Range("D4").Value = 0.0987583780731366
$ws.Range("E4").Value = 0.01263758662861791

$ws.Range("D5").Value = 0.1015083522658633
$ws.Range("E5").Value = -0.001733102253032803

$ws.Range("D6").Value = 0.05940578140068804
$ws.Range("E6").Value = -0.006624029237094486

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.00615048046081168

# Restore sheet protection (matching original protected state)
$ws.Protect()
